$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Treatment Agent" query (B5 / TreatmentTab row): drop the redundant
# CONCAT() wrapper around REPLACE(...), per commit "Fixed query issue for
# C3DC phs002599".
$fixedTreatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs002599' AND srv.first_event IN ('Not Reported')
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

$ws.Cells.Item(5, 2).Value = $fixedTreatmentQuery

# Re-apply the (unchanged) font so the engine mints a fresh style record for
# the edited cell, matching how Excel re-stamped B5's style on save.
$ws.Cells.Item(5, 2).Font.Name = $ws.Cells.Item(5, 2).Font.Name

# The author's edit also left the sheet scrolled/selected one row higher
# (A6/C7 -> A5/B5), i.e. focused back on the just-fixed Treatment cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("B5").Select()
